# Updated titles for iPortal users
# Replace the placeholder "TBD" title values (rows 35-44, column C) on the
# "Attributes" sheet with the real job titles for each user.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")
$ws.Activate()

$ws.Range("C35").Value = "BI Administrator"
$ws.Range("C36").Value = "BI Developer"
$ws.Range("C37").Value = "SR. BI Developr"
$ws.Range("C38").Value = "BI Analyst"
$ws.Range("C39").Value = "Data Analyst"
$ws.Range("C40").Value = "Marketing Analyst"
$ws.Range("C41").Value = "Sales Operations Lead"
$ws.Range("C42").Value = "Financial Analyst"
$ws.Range("C43").Value = "VP - Operations"
$ws.Range("C44").Value = "Sales Territory Mgr"

# Match the author's on-screen view state after making the edits.
$ws.Range("C45").Select()
